$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / Row 11: simple value bumps
$ws.Range("B2").Value = 91829
$ws.Range("B11").Value = 91809

# Rows 12 and 13 swap their species data (with B13 getting a new value
# rather than the plain swapped one).

# Target row 12 (previously row 13's data)
$ws.Range("A12").Value = 131130472
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("Q12").Value = 410737
$ws.Range("R12").Value = 7037762
$ws.Range("Z12").Value = "15:12"
$ws.Range("AB12").Value = "15:12"

# Target row 13 (previously row 12's data)
$ws.Range("A13").Value = 131131304
$ws.Range("B13").Value = 91829
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = "Granticka"
$ws.Range("G13").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H13").Value = ""
$ws.Range("Q13").Value = 410603
$ws.Range("R13").Value = 7037541
$ws.Range("Z13").Value = "16:00"
$ws.Range("AB13").Value = "16:00"
